$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '23.092.77'
$ws.Range("E2").Value = '  -3.39%  '
$ws.Range("D3").Value = '1.604.59'
$ws.Range("E3").Value = '  -2.87%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("D6").Value = '301.61'
$ws.Range("E6").Value = '  -3.02%  '
$ws.Range("D7").Value = '0.3784'
$ws.Range("E7").Value = '  -2.85%  '
$ws.Range("D8").Value = '0.3658'
$ws.Range("E8").Value = '  -4.35%  '
$ws.Range("D9").Value = '50.00'
$ws.Range("E9").Value = '  -3.64%  '
$ws.Range("D10").Value = '1.269'
$ws.Range("E10").Value = '  -5.48%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '0.08153'
$ws.Range("E11").Value = '  -3.49%  '
$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("D13").Value = '22.88'
$ws.Range("E13").Value = '  -4.39%  '
$ws.Range("D14").Value = '6.612'
$ws.Range("E14").Value = '  -5.83%  '
$ws.Range("D15").Value = '0.00001258'
$ws.Range("E15").Value = '  -4.28%  '
$ws.Range("D16").Value = '7.397'
$ws.Range("E16").Value = '  -7.83%  '
$ws.Range("D17").Value = '1.599.23'
$ws.Range("E17").Value = '  -3.08%  '
$ws.Range("D18").Value = '92.12'
$ws.Range("E18").Value = '  -2.14%  '
$ws.Range("D19").Value = '0.06888'
$ws.Range("E19").Value = '  -1.37%  '
$ws.Range("D20").Value = '18.30'
$ws.Range("E20").Value = '  -6.59%  '
$ws.Range("D21").Value = '6.589'
$ws.Range("E21").Value = '  -5.30%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").Value = '13.13'
$ws.Range("E23").Value = '  -4.00%  '
$ws.Range("B24").Value = 'WrappedBTC'
$ws.Range("C24").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D24").Value = '23.094.11'
$ws.Range("E24").Value = '  -3.33%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '2.353'
$ws.Range("E25").Value = '  -3.67%  '
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = '2.799'
$ws.Range("E26").Value = '  -4.95%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '21.17'
$ws.Range("E27").Value = '  -3.87%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = '150.52'
$ws.Range("E28").Value = '  -1.72%  '
$ws.Range("B29").Value = 'HuobiToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D29").Value = '5.274'
$ws.Range("E29").Value = '  -2.44%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = '133.73'
$ws.Range("E30").Value = '  -2.66%  '
$ws.Range("B31").Value = 'WEMIXTOKEN'
$ws.Range("C31").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D31").Value = '2.370'
$ws.Range("E31").Value = '  -4.57%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '6.847'
$ws.Range("E32").Value = '  -11.56%  '
$ws.Range("B33").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C33").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D33").Value = '1.774.57'
$ws.Range("E33").Value = '  -3.09%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '0.9557'
$ws.Range("E34").Value = '  -4.15%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.07701'
$ws.Range("E35").Value = '  -5.27%  '
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").Value = '10.46'
$ws.Range("E36").Value = '  -2.33%  '
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '6.292'
$ws.Range("E37").Value = '  -5.71%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.02725'
$ws.Range("E38").Value = '  -6.07%  '
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").Value = '0.2548'
$ws.Range("E39").Value = '  -4.81%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '0.08920'
$ws.Range("E40").Value = '  -1.97%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '1.369'
$ws.Range("E41").Value = '  -3.57%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.7094'
$ws.Range("E42").Value = '  -6.19%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '12.69'
$ws.Range("E43").Value = '  -6.13%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '15.30'
$ws.Range("E44").Value = '  -6.56%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '0.6641'
$ws.Range("E45").Value = '  -4.39%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = '2.336'
$ws.Range("E46").Value = '  -4.29%  '
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  +0.24%  '
$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").Value = '4.003'
$ws.Range("E48").Value = '  -2.63%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '132.62'
$ws.Range("E49").Value = '  -0.62%  '
$ws.Range("B50").Value = 'Flow'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D50").Value = '1.246'
$ws.Range("E50").Value = '  +1.71%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.07936'
$ws.Range("E51").Value = '  -4.02%  '
